$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Strip the "Heading2" paragraph style from the seven section headings,
#    turning them back into plain (Normal) paragraphs.
# ---------------------------------------------------------------------------
$headingIndexes = 2,6,12,18,26,36,42
foreach ($idx in $headingIndexes) {
    $p = $d.Paragraphs($idx)
    $p.Style = "Normal"
}

# ---------------------------------------------------------------------------
# 2) Replace the inline citation markers with the new "Ref-XXXXXXX" /
#    author-style tokens. Each replacement is scoped to the specific
#    paragraph's Range (rather than the whole document) because the same
#    source citation text (e.g. "(Datar and Rajan)", "(Javaid et al.)")
#    maps to different replacement tokens in different paragraphs.
# ---------------------------------------------------------------------------
function Replace-InParagraph($paraIndex, $oldText, $newText) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $r.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
}

# Paragraph 8 - "Current Strategy Analysis of C&KM", first body paragraph
Replace-InParagraph 8 "(Datar and Rajan)" "(Ref-f361053)"
Replace-InParagraph 8 "(Stone et al.)" "(Ref-f361053)"

# Paragraph 10 - "Furthermore, C&KM's strategic choice..."
Replace-InParagraph 10 "(Stone et al.)" "(Ref-s606200)"

# Paragraph 14 - "Stone Manufacturing exemplifies..."
Replace-InParagraph 14 "(Rai et al.)" "(Ref-u418791)"
Replace-InParagraph 14 "(Javaid et al.)" "(Ref-u418791)"

# Paragraph 16 - "Additionally, Stone Manufacturing's strategic emphasis..."
Replace-InParagraph 16 "(Javaid et al.)" "(Ref-s825060)"

# Paragraph 20 - "To enhance product quality and reduce costs..."
Replace-InParagraph 20 "(He and Bai)" "(Ref-u491647)"
Replace-InParagraph 20 "(Rai et al.)" "(Ref-u491647)"

# Paragraph 22 - "Moreover, process improvements at C&KM..."
Replace-InParagraph 22 "(He and Bai)" "(Ref-f250134)"

# Paragraph 24 - "In addition, worker training is pivotal..."
Replace-InParagraph 24 "(Rai et al.)" "(Ref-f375515)"

# Paragraph 28 - "To effectively address the financial perspective..."
Replace-InParagraph 28 "(Frederico et al.)" "(Brown and Garcia)"

# Paragraph 30 - "Similarly, the customer perspective..."
Replace-InParagraph 30 "(Frederico et al.)" "(Ref-s373905)"

# Paragraph 32 - "Likewise, the internal business processes perspective..."
Replace-InParagraph 32 "(He and Bai)" "(Ref-u289928)"
Replace-InParagraph 32 "(Frederico et al.)" "(Ref-u289928)"

# Paragraph 34 - "Furthermore, establishing criteria for the learning..."
Replace-InParagraph 34 "(Frederico et al.)" "(Ref-s898831)"

# Paragraph 38 - "Implementing the proposed strategies..."
Replace-InParagraph 38 "(Datar and Rajan)" "(Ref-f732504)"

# Paragraph 40 - "Therefore, continuous evaluation and adaptation..."
Replace-InParagraph 40 "(Javaid et al.)" "(Johnson 45)"

Write-Output "edit complete"
